$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26, 27, 28 (columns A, B, E, F, G, H, Q, R) are cyclically rotated:
#   new row26 = old row28
#   new row27 = old row26
#   new row28 = old row27
# Capture the "before" values first (via Value2, which returns the raw
# scalar rather than a property descriptor), then write them into place.

$vA26 = $ws.Range("A26").Value2
$vB26 = $ws.Range("B26").Value2
$vE26 = $ws.Range("E26").Value2
$vF26 = $ws.Range("F26").Value2
$vG26 = $ws.Range("G26").Value2
$vH26 = $ws.Range("H26").Value2
$vQ26 = $ws.Range("Q26").Value2
$vR26 = $ws.Range("R26").Value2

$vA27 = $ws.Range("A27").Value2
$vB27 = $ws.Range("B27").Value2
$vE27 = $ws.Range("E27").Value2
$vF27 = $ws.Range("F27").Value2
$vG27 = $ws.Range("G27").Value2
$vH27 = $ws.Range("H27").Value2
$vQ27 = $ws.Range("Q27").Value2
$vR27 = $ws.Range("R27").Value2

$vA28 = $ws.Range("A28").Value2
$vB28 = $ws.Range("B28").Value2
$vE28 = $ws.Range("E28").Value2
$vF28 = $ws.Range("F28").Value2
$vG28 = $ws.Range("G28").Value2
$vH28 = $ws.Range("H28").Value2
$vQ28 = $ws.Range("Q28").Value2
$vR28 = $ws.Range("R28").Value2

# new row26 = old row28
$ws.Range("A26").Value = $vA28
$ws.Range("B26").Value = $vB28
$ws.Range("E26").Value = $vE28
$ws.Range("F26").Value = $vF28
$ws.Range("G26").Value = $vG28
$ws.Range("H26").Value = $vH28
$ws.Range("Q26").Value = $vQ28
$ws.Range("R26").Value = $vR28

# new row27 = old row26
$ws.Range("A27").Value = $vA26
$ws.Range("B27").Value = $vB26
$ws.Range("E27").Value = $vE26
$ws.Range("F27").Value = $vF26
$ws.Range("G27").Value = $vG26
$ws.Range("H27").Value = $vH26
$ws.Range("Q27").Value = $vQ26
$ws.Range("R27").Value = $vR26

# new row28 = old row27
$ws.Range("A28").Value = $vA27
$ws.Range("B28").Value = $vB27
$ws.Range("E28").Value = $vE27
$ws.Range("F28").Value = $vF27
$ws.Range("G28").Value = $vG27
$ws.Range("H28").Value = $vH27
$ws.Range("Q28").Value = $vQ27
$ws.Range("R28").Value = $vR27
